$d = $word.ActiveDocument

# 1. Merge the split runs around "executes." (this also drops the stray
#    mid-sentence _GoBack bookmark that used to separate "exe" / "cutes.")
$d.Content.Find.Execute(
    "state while it executes. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "state while it executes. ", 2) | Out-Null

# 2. Append the new sentence to the "nodes are simple python functions..."
#    paragraph (paragraph containing that sentence).
$nodesFuncPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*simple python functions*") {
        $nodesFuncPara = $p
        break
    }
}
$nodesFuncPara.Range.InsertAfter(" So, It can be determined that nodes take state as an input and return state as the output.")

# 3. Insert a brand-new list paragraph right after it, with the doc-string
#    guidance text, and put the (hidden) _GoBack bookmark at the very end
#    of that new paragraph's content.
$nodesFuncPara.Range.InsertParagraphAfter()
$newPara = $nodesFuncPara.Next()

$docStringText = "While building nodes, emphasis should be on providing doc strings to each and every node. Because eventually, we will be building AI agents. For these to have context of what each node is doing, Doc Strings are necessary."

# Write the text with one extra sentinel character so the bookmark we need
# to plant sits *inside* the run (not exactly at the paragraph's text end,
# which the COM bookmark placement mishandles for the reserved "_GoBack"
# name). We then trim the sentinel back off.
$newPara.Range.Text = $docStringText + "#"

$bmStart = $newPara.Range.Start + $docStringText.Length
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($bmStart, $bmStart + 1)
$sentinelRange.Delete()

foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [$($p.Range.Text)]"
}
